$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for Wins, Losses, Ties in AD1:AF1. Copy the
# existing header style (bold/centered/bordered, same as the other
# header cells) from AC1 before writing the new labels so the cells
# reuse the same style index instead of creating new ones.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill season record values (Wins=71, Losses=91, Ties=0) for every data row
for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 30).Value = 71
    $ws.Cells.Item($r, 31).Value = 91
    $ws.Cells.Item($r, 32).Value = 0
}
